$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"

# Column D keeps the same date-time number format as the rows above it.
$ws.Cells.Item($row, 4).Value = 44448
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112052
$ws.Cells.Item($row, 7).Value = "Albahaca"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 640
$ws.Cells.Item($row, 11).Value = 4500
$ws.Cells.Item($row, 12).Value = 5000
$ws.Cells.Item($row, 13).Value = 4750
$ws.Cells.Item($row, 14).Value = "$/paquete"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 4750
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
